$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting old row 2 down to row 3
$ws.Rows("2:2").Insert()

# New row 2: identity_type (F) and time_type (J)
$ws.Cells.Item(2, 6).Value = 1

$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "12:13:14.987654321"

# Update shifted row 3: identity_type value incremented
$ws.Cells.Item(3, 6).Value = 2

# Update shifted row 3: uuid_type value changed
$ws.Cells.Item(3, 19).Value = "6779defb-6d49-4e2e-b3dd-95cd071cea5c"
